$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-01-31 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-02-01 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("23+8=", $true, $false, $false, $false, $false, $true, 1, $false, "22+22=", 2) | Out-Null
$d.Content.Find.Execute("24+41=", $true, $false, $false, $false, $false, $true, 1, $false, "9+1=", 2) | Out-Null
$d.Content.Find.Execute("12-11=", $true, $false, $false, $false, $false, $true, 1, $false, "20-12=", 2) | Out-Null
$d.Content.Find.Execute("51-18=", $true, $false, $false, $false, $false, $true, 1, $false, "0+45=", 2) | Out-Null
$d.Content.Find.Execute("79-21=", $true, $false, $false, $false, $false, $true, 1, $false, "32-26=", 2) | Out-Null
$d.Content.Find.Execute("97-81=", $true, $false, $false, $false, $false, $true, 1, $false, "6+73=", 2) | Out-Null
$d.Content.Find.Execute("96-37=", $true, $false, $false, $false, $false, $true, 1, $false, "58+21=", 2) | Out-Null
$d.Content.Find.Execute("98-24=", $true, $false, $false, $false, $false, $true, 1, $false, "42-10=", 2) | Out-Null
$d.Content.Find.Execute("4+71=", $true, $false, $false, $false, $false, $true, 1, $false, "59+36=", 2) | Out-Null
$d.Content.Find.Execute("77+18=", $true, $false, $false, $false, $false, $true, 1, $false, "18+1=", 2) | Out-Null
$d.Content.Find.Execute("62-1=", $true, $false, $false, $false, $false, $true, 1, $false, "82+15=", 2) | Out-Null
$d.Content.Find.Execute("86-6=", $true, $false, $false, $false, $false, $true, 1, $false, "96-28=", 2) | Out-Null
$d.Content.Find.Execute("81+3=", $true, $false, $false, $false, $false, $true, 1, $false, "12+30=", 2) | Out-Null
$d.Content.Find.Execute("46-36=", $true, $false, $false, $false, $false, $true, 1, $false, "30-12=", 2) | Out-Null
$d.Content.Find.Execute("34-31=", $true, $false, $false, $false, $false, $true, 1, $false, "26+7=", 2) | Out-Null
$d.Content.Find.Execute("73-26=", $true, $false, $false, $false, $false, $true, 1, $false, "28+49=", 2) | Out-Null
$d.Content.Find.Execute("9+44=", $true, $false, $false, $false, $false, $true, 1, $false, "2+53=", 2) | Out-Null
$d.Content.Find.Execute("24+62=", $true, $false, $false, $false, $false, $true, 1, $false, "73-50=", 2) | Out-Null
$d.Content.Find.Execute("32+17=", $true, $false, $false, $false, $false, $true, 1, $false, "40+38=", 2) | Out-Null
$d.Content.Find.Execute("31+8=", $true, $false, $false, $false, $false, $true, 1, $false, "46+0=", 2) | Out-Null
$d.Content.Find.Execute("74-63=", $true, $false, $false, $false, $false, $true, 1, $false, "56-5=", 2) | Out-Null
$d.Content.Find.Execute("11+34=", $true, $false, $false, $false, $false, $true, 1, $false, "79-24=", 2) | Out-Null
$d.Content.Find.Execute("5+53=", $true, $false, $false, $false, $false, $true, 1, $false, "90-20=", 2) | Out-Null
$d.Content.Find.Execute("72+0=", $true, $false, $false, $false, $false, $true, 1, $false, "48-43=", 2) | Out-Null
$d.Content.Find.Execute("90-63=", $true, $false, $false, $false, $false, $true, 1, $false, "61-21=", 2) | Out-Null
$d.Content.Find.Execute("22+53=", $true, $false, $false, $false, $false, $true, 1, $false, "65-9=", 2) | Out-Null
$d.Content.Find.Execute("42+9=", $true, $false, $false, $false, $false, $true, 1, $false, "99-3=", 2) | Out-Null
$d.Content.Find.Execute("82-82=", $true, $false, $false, $false, $false, $true, 1, $false, "14+65=", 2) | Out-Null
$d.Content.Find.Execute("24+35=", $true, $false, $false, $false, $false, $true, 1, $false, "72-68=", 2) | Out-Null
$d.Content.Find.Execute("13+9=", $true, $false, $false, $false, $false, $true, 1, $false, "88-37=", 2) | Out-Null
$d.Content.Find.Execute("50-12=", $true, $false, $false, $false, $false, $true, 1, $false, "37-17=", 2) | Out-Null
$d.Content.Find.Execute("90+8=", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=", 2) | Out-Null
$d.Content.Find.Execute("7+52=", $true, $false, $false, $false, $false, $true, 1, $false, "81-80=", 2) | Out-Null
$d.Content.Find.Execute("90-78=", $true, $false, $false, $false, $false, $true, 1, $false, "66+19=", 2) | Out-Null
$d.Content.Find.Execute("91-90=", $true, $false, $false, $false, $false, $true, 1, $false, "8+14=", 2) | Out-Null
$d.Content.Find.Execute("69-14=", $true, $false, $false, $false, $false, $true, 1, $false, "45+43=", 2) | Out-Null
$d.Content.Find.Execute("48+51=", $true, $false, $false, $false, $false, $true, 1, $false, "72-7=", 2) | Out-Null
$d.Content.Find.Execute("10+63=", $true, $false, $false, $false, $false, $true, 1, $false, "52-50=", 2) | Out-Null
$d.Content.Find.Execute("40-24=", $true, $false, $false, $false, $false, $true, 1, $false, "13+72=", 2) | Out-Null
$d.Content.Find.Execute("85+6=", $true, $false, $false, $false, $false, $true, 1, $false, "7+43=", 2) | Out-Null
$d.Content.Find.Execute("88-23=", $true, $false, $false, $false, $false, $true, 1, $false, "25+70=", 2) | Out-Null
$d.Content.Find.Execute("84-28=", $true, $false, $false, $false, $false, $true, 1, $false, "50-1=", 2) | Out-Null
$d.Content.Find.Execute("27+43=", $true, $false, $false, $false, $false, $true, 1, $false, "37+27=", 2) | Out-Null
$d.Content.Find.Execute("91-54=", $true, $false, $false, $false, $false, $true, 1, $false, "67-52=", 2) | Out-Null
$d.Content.Find.Execute("94+4=", $true, $false, $false, $false, $false, $true, 1, $false, "49-20=", 2) | Out-Null
$d.Content.Find.Execute("64-24=", $true, $false, $false, $false, $false, $true, 1, $false, "48+41=", 2) | Out-Null
$d.Content.Find.Execute("31+32=", $true, $false, $false, $false, $false, $true, 1, $false, "42-3=", 2) | Out-Null
$d.Content.Find.Execute("31+13=", $true, $false, $false, $false, $false, $true, 1, $false, "87-45=", 2) | Out-Null
$d.Content.Find.Execute("42-20=", $true, $false, $false, $false, $false, $true, 1, $false, "20+5=", 2) | Out-Null
$d.Content.Find.Execute("60-53=", $true, $false, $false, $false, $false, $true, 1, $false, "22+19=", 2) | Out-Null
$d.Content.Find.Execute("20+61=", $true, $false, $false, $false, $false, $true, 1, $false, "25+57=", 2) | Out-Null
$d.Content.Find.Execute("67-46=", $true, $false, $false, $false, $false, $true, 1, $false, "71-37=", 2) | Out-Null
$d.Content.Find.Execute("0+47=", $true, $false, $false, $false, $false, $true, 1, $false, "62-9=", 2) | Out-Null
$d.Content.Find.Execute("29+61=", $true, $false, $false, $false, $false, $true, 1, $false, "32+64=", 2) | Out-Null
$d.Content.Find.Execute("13+22=", $true, $false, $false, $false, $false, $true, 1, $false, "58-1=", 2) | Out-Null
$d.Content.Find.Execute("84-50=", $true, $false, $false, $false, $false, $true, 1, $false, "64+10=", 2) | Out-Null
$d.Content.Find.Execute("19+26=", $true, $false, $false, $false, $false, $true, 1, $false, "10+20=", 2) | Out-Null
$d.Content.Find.Execute("50-31=", $true, $false, $false, $false, $false, $true, 1, $false, "99-16=", 2) | Out-Null
$d.Content.Find.Execute("55+2=", $true, $false, $false, $false, $false, $true, 1, $false, "64-46=", 2) | Out-Null
$d.Content.Find.Execute("43-8=", $true, $false, $false, $false, $false, $true, 1, $false, "59-44=", 2) | Out-Null
$d.Content.Find.Execute("36+45=", $true, $false, $false, $false, $false, $true, 1, $false, "69-17=", 2) | Out-Null
$d.Content.Find.Execute("63-51=", $true, $false, $false, $false, $false, $true, 1, $false, "10-3=", 2) | Out-Null
$d.Content.Find.Execute("56-23=", $true, $false, $false, $false, $false, $true, 1, $false, "5+74=", 2) | Out-Null
$d.Content.Find.Execute("72-27=", $true, $false, $false, $false, $false, $true, 1, $false, "35-15=", 2) | Out-Null
$d.Content.Find.Execute("15+4=", $true, $false, $false, $false, $false, $true, 1, $false, "7+18=", 2) | Out-Null
$d.Content.Find.Execute("91+0=", $true, $false, $false, $false, $false, $true, 1, $false, "37-15=", 2) | Out-Null
$d.Content.Find.Execute("30-9=", $true, $false, $false, $false, $false, $true, 1, $false, "57-19=", 2) | Out-Null
$d.Content.Find.Execute("41-33=", $true, $false, $false, $false, $false, $true, 1, $false, "18+32=", 2) | Out-Null
$d.Content.Find.Execute("90-69=", $true, $false, $false, $false, $false, $true, 1, $false, "45-9=", 2) | Out-Null
$d.Content.Find.Execute("41+43=", $true, $false, $false, $false, $false, $true, 1, $false, "27+48=", 2) | Out-Null
$d.Content.Find.Execute("40+21=", $true, $false, $false, $false, $false, $true, 1, $false, "97-76=", 2) | Out-Null
$d.Content.Find.Execute("46-19=", $true, $false, $false, $false, $false, $true, 1, $false, "48-26=", 2) | Out-Null
$d.Content.Find.Execute("34+30=", $true, $false, $false, $false, $false, $true, 1, $false, "13+11=", 2) | Out-Null
$d.Content.Find.Execute("52+17=", $true, $false, $false, $false, $false, $true, 1, $false, "62-19=", 2) | Out-Null
$d.Content.Find.Execute("68+12=", $true, $false, $false, $false, $false, $true, 1, $false, "37+52=", 2) | Out-Null
$d.Content.Find.Execute("57+13=", $true, $false, $false, $false, $false, $true, 1, $false, "12+38=", 2) | Out-Null
$d.Content.Find.Execute("95-90=", $true, $false, $false, $false, $false, $true, 1, $false, "78+6=", 2) | Out-Null
$d.Content.Find.Execute("5+62=", $true, $false, $false, $false, $false, $true, 1, $false, "54+3=", 2) | Out-Null
$d.Content.Find.Execute("41+39=", $true, $false, $false, $false, $false, $true, 1, $false, "65+31=", 2) | Out-Null
$d.Content.Find.Execute("7+50=", $true, $false, $false, $false, $false, $true, 1, $false, "64+32=", 2) | Out-Null
$d.Content.Find.Execute("96-30=", $true, $false, $false, $false, $false, $true, 1, $false, "25-12=", 2) | Out-Null
$d.Content.Find.Execute("76-40=", $true, $false, $false, $false, $false, $true, 1, $false, "53+26=", 2) | Out-Null
$d.Content.Find.Execute("66-43=", $true, $false, $false, $false, $false, $true, 1, $false, "76-2=", 2) | Out-Null
$d.Content.Find.Execute("75+24=", $true, $false, $false, $false, $false, $true, 1, $false, "35+3=", 2) | Out-Null
$d.Content.Find.Execute("72-4=", $true, $false, $false, $false, $false, $true, 1, $false, "29+10=", 2) | Out-Null
$d.Content.Find.Execute("87-0=", $true, $false, $false, $false, $false, $true, 1, $false, "33+22=", 2) | Out-Null
$d.Content.Find.Execute("23+23=", $true, $false, $false, $false, $false, $true, 1, $false, "81-79=", 2) | Out-Null
$d.Content.Find.Execute("22+8=", $true, $false, $false, $false, $false, $true, 1, $false, "26+65=", 2) | Out-Null
$d.Content.Find.Execute("47+46=", $true, $false, $false, $false, $false, $true, 1, $false, "61+6=", 2) | Out-Null
$d.Content.Find.Execute("33+59=", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=", 2) | Out-Null
$d.Content.Find.Execute("28+18=", $true, $false, $false, $false, $false, $true, 1, $false, "39-13=", 2) | Out-Null
$d.Content.Find.Execute("7+11=", $true, $false, $false, $false, $false, $true, 1, $false, "49-28=", 2) | Out-Null
$d.Content.Find.Execute("75-27=", $true, $false, $false, $false, $false, $true, 1, $false, "13+50=", 2) | Out-Null
$d.Content.Find.Execute("88-72=", $true, $false, $false, $false, $false, $true, 1, $false, "50-47=", 2) | Out-Null
$d.Content.Find.Execute("73-6=", $true, $false, $false, $false, $false, $true, 1, $false, "64+13=", 2) | Out-Null
$d.Content.Find.Execute("20-18=", $true, $false, $false, $false, $false, $true, 1, $false, "68+4=", 2) | Out-Null
$d.Content.Find.Execute("42+54=", $true, $false, $false, $false, $false, $true, 1, $false, "98-20=", 2) | Out-Null
$d.Content.Find.Execute("9+17=", $true, $false, $false, $false, $false, $true, 1, $false, "1+83=", 2) | Out-Null
$d.Content.Find.Execute("23+53=", $true, $false, $false, $false, $false, $true, 1, $false, "34+15=", 2) | Out-Null
$d.Content.Find.Execute("7+32=", $true, $false, $false, $false, $false, $true, 1, $false, "31+10=", 2) | Out-Null
